$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shift existing rows to make room for new "account / password" translation
# rows, and leave a blank separator row, matching the shape of the target
# sheet (old row 3 "lang_register" ends up at row 5; old rows 4-10 end up at
# rows 12-18; old rows 11-19 end up at rows 19-27). Row 11 stays blank.
# ---------------------------------------------------------------------------

# Make room for 2 new rows (future rows 3,4) before old row 3
$ws.Rows.Item(3).Resize(2).Insert()

# Old row 3 (lang_register) is now row 5. Make room for 5 new rows
# (future rows 6-10) before old row 4 (now at row 6)
$ws.Rows.Item(6).Resize(5).Insert()

# Old row 4 (lang_click) is now at row 11. Insert one blank separator row
# before it so old row 4 becomes row 12 and a blank row 11 remains.
$ws.Rows.Item(11).Resize(1).Insert()

# ---------------------------------------------------------------------------
# Fill the newly-inserted rows with the new login/account/password strings
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "lang_change_password"
$ws.Range("B3").Value = "Đổi Mật Khẩu"
$ws.Range("C3").Value = "Change Password"

$ws.Range("A4").Value = "lang_forget_password"
$ws.Range("B4").Value = "Chưa Có/ Quên Mật Khẩu"
$ws.Range("C4").Value = "Forget Password"

$ws.Range("A6").Value = "lang_phone_number"
$ws.Range("B6").Value = "Số điện thoại"
$ws.Range("C6").Value = "Phone number"

$ws.Range("A7").Value = "lang_email_receive_password"
$ws.Range("B7").Value = "Địa chỉ email nhận mật khẩu"
$ws.Range("C7").Value = "Address receive password"

$ws.Range("A8").Value = "lang_student_code"
$ws.Range("B8").Value = "Mã số sinh viên"
$ws.Range("C8").Value = "Student ID"

$ws.Range("A9").Value = "lang_send_get_login_info"
$ws.Range("B9").Value = "Gửi Thông Tin Xác Thực"
$ws.Range("C9").Value = "Request Log In Info"

$ws.Range("A10").Value = "lang_back_to_login"
$ws.Range("B10").Value = "Quay Về Trang Đăng Nhập"
$ws.Range("C10").Value = "Back To LogIn Page"

# ---------------------------------------------------------------------------
# Append two more rows (28, 29) at the end, matching the style (font/valign)
# of the preceding data rows (19-27), by copying row 27's formatting down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Copy()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Copy()
$ws.Rows.Item(29).Insert()

$ws.Range("A28").Value = "lang_account"
$ws.Range("B28").Value = "Tài khoản"
$ws.Range("C28").Value = "Account"

$ws.Range("A29").Value = "lang_password"
$ws.Range("B29").Value = "Mật khẩu"
$ws.Range("C29").Value = "Password"

# ---------------------------------------------------------------------------
# Column widths: A stays the same; B and C grow to fit the new, longer text
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 41.3
$ws.Columns.Item(3).ColumnWidth = 31.8

# ---------------------------------------------------------------------------
# Selection moves from A11 to A10
# ---------------------------------------------------------------------------
$ws.Range("A10").Select()
